$wb = $excel.ActiveWorkbook

# "展览" sheet: row14 (F14) and row16 (F16) "想去人数" counts updated
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F14").Value = 1377
$wsExhibit.Range("F16").Value = 777

# "全部类型" sheet: same rows, but F16 updated to a different value
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 1377
$wsAll.Range("F16").Value = 778
